$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.9999674344715328,
    0.9989400190852257,
    0.9999604067229327,
    0.9999714913392064,
    0.9999652011360065,
    0.00003039849136261971,
    0.0009894456561571146,
    0.00003825735474157083,
    0.00002064418294921036,
    0.0000294507688453906,
    0.000349036811621329,
    0.005513482689065026,
    1.000060120975632,
    0.005748203061621094,
    94.80223515402628,
    139.9006406741497
)

for ($row = 2; $row -le 26; $row++) {
    for ($i = 0; $i -lt $newValues.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $newValues[$i]
    }
}
